$d = $word.ActiveDocument

# Locate the paragraph that contains exactly "dotnet build" (the one that
# precedes the "Static code analysis" section which is being removed).
$targetIndex = -1
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    if ($p.Range.Text.Trim() -eq "dotnet build") {
        $targetIndex = $i
    }
}

# The paragraph right after "dotnet build" is the (currently empty)
# paragraph that we will keep and reuse as the new trailing paragraph.
$keepPara = $d.Paragraphs.Item($targetIndex + 1)

# Delete everything from the end of that kept paragraph through to the
# very end of the document body (this removes the whole "Static code
# analysis" / StyleCop section, including the old _GoBack bookmark that
# lived further down, and the trailing empty paragraphs).
$delStart = $keepPara.Range.End
$delEnd = $d.Content.End
if ($delEnd -gt $delStart) {
    $delRange = $d.Range($delStart, $delEnd)
    $delRange.Delete()
}

# Re-create the "_GoBack" bookmark as an empty (collapsed) bookmark inside
# the kept trailing paragraph. Adding a bookmark directly on a zero-length
# range inside an empty paragraph doesn't anchor correctly, so as a
# workaround we insert a couple of placeholder characters, wrap the
# bookmark around them, then delete the placeholder text - leaving the
# bookmark collapsed in place.
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$r = $lastPara.Range
$r.InsertAfter("ZZ")
$bmRange = $d.Range($r.Start, $r.Start + 2)
$d.Bookmarks.Add("_GoBack", $bmRange)
$placeholder = $d.Range($r.Start, $r.Start + 2)
$placeholder.Text = ""
